$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-02-23 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-02-24 Saturday", 2)

# Update the division-problem table cells. The table has 20 rows x 5 cols,
# but only every 4th row (1, 5, 9, 13, 17) actually contains text.
$table = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)

$oldValues = @(
    @("43÷3=14, 1", "37÷9=4, 1", "91÷5=18, 1", "24÷2=12, 0", "40÷2=20, 0"),
    @("27÷6=4, 3", "19÷8=2, 3", "27÷8=3, 3", "25÷3=8, 1", "45÷4=11, 1"),
    @("22÷3=7, 1", "20÷7=2, 6", "59÷8=7, 3", "39÷2=19, 1", "90÷9=10, 0"),
    @("68÷7=9, 5", "89÷2=44, 1", "82÷2=41, 0", "79÷3=26, 1", "61÷8=7, 5"),
    @("68÷7=9, 5", "56÷4=14, 0", "80÷4=20, 0", "23÷5=4, 3", "91÷9=10, 1")
)

$newValues = @(
    @("53÷8=6, 5", "94÷7=13, 3", "99÷8=12, 3", "83÷8=10, 3", "21÷3=7, 0"),
    @("83÷7=11, 6", "42÷7=6, 0", "95÷8=11, 7", "36÷2=18, 0", "62÷3=20, 2"),
    @("93÷3=31, 0", "58÷4=14, 2", "80÷3=26, 2", "77÷2=38, 1", "25÷5=5, 0"),
    @("61÷4=15, 1", "40÷3=13, 1", "42÷8=5, 2", "61÷3=20, 1", "57÷8=7, 1"),
    @("48÷5=9, 3", "22÷7=3, 1", "60÷7=8, 4", "36÷3=12, 0", "78÷3=26, 0")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowIndex = $rows[$r]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $table.Cell($rowIndex, $c)
        $cellRange = $cell.Range
        $cellRange.MoveEnd(1, -1) | Out-Null
        $old = $oldValues[$r][$c - 1]
        $new = $newValues[$r][$c - 1]
        if ($cellRange.Text -eq $old) {
            $cellRange.Text = $new
        } else {
            $cellRange.Find.Execute($old, $true, $false, $false, $false, $false, `
                                     $true, 1, $false, $new, 2)
        }
    }
}
